# Update countries & provincias Spain
#
# 1) Shared-strings reorder: "Niger" moves from right after "Libano" to
#    right before "Principado de Andorra" (i.e. right after "Kirguistan").
#    Column A of the data rows keeps referencing the same shared-string
#    slot, so moving the text itself re-labels rows 99-101 (which used to
#    read Principado de Andorra / Costa Rica / Libano) while row 98 (which
#    used to read Principado de Andorra) now reads "Niger".
#
# 2) Because the row 98 label is now "Niger", its statistics (columns B-H)
#    are replaced with fresh Niger numbers, and the old Principado de
#    Andorra / Costa Rica / Libano statistics shift down one row each
#    (99<-98, 100<-99, 101<-100). Row 102 (Somalia) is untouched.
#
# 3) A handful of unrelated totals were refreshed for Estados Unidos (row 4),
#    Alemania (row 9) and Colombia (row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-order the "Niger" entry in the shared country list -------------
$ws.Range("A98").Value = "Niger"
$ws.Range("A99").Value = "Principado de Andorra"
$ws.Range("A100").Value = "Costa Rica"
$ws.Range("A101").Value = "Libano"

# --- 2) Shift the Principado de Andorra / Costa Rica / Libano stats down,
#        and give the newly relabeled Niger row (98) its own fresh figures.

# Row 101 <- old row 100 values (Libano)
$ws.Range("B101").Value = 737
$ws.Range("C101").Value = 4
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = 512
$ws.Range("F101").Value = 43
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 25

# Row 100 <- old row 99 values (Costa Rica)
$ws.Range("B100").Value = 739
$ws.Range("C100").Value = 6
$ws.Range("D100").Value = 386
$ws.Range("E100").Value = 347
$ws.Range("F100").Value = 6
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 6

# Row 99 <- old row 98 values (Principado de Andorra)
$ws.Range("B99").Value = 747
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 472
$ws.Range("E99").Value = 231
$ws.Range("F99").Value = 17
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 44

# Row 98 <- brand new Niger values
$ws.Range("B98").Value = 750
$ws.Range("C98").Value = 14
$ws.Range("D98").Value = 518
$ws.Range("E98").Value = 196
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 36

# --- 3) Refresh Estados Unidos (row 4) -------------------------------------
$ws.Range("B4").Value = 1185285
$ws.Range("C4").Value = 24511
$ws.Range("E4").Value = 938559
$ws.Range("G4").Value = 1063
$ws.Range("H4").Value = 68507

# --- Refresh Alemania (row 9) ----------------------------------------------
$ws.Range("B9").Value = 165664
$ws.Range("C9").Value = 697
$ws.Range("E9").Value = 28198
$ws.Range("G9").Value = 54
$ws.Range("H9").Value = 6866

# --- Refresh Colombia (row 48) ----------------------------------------------
$ws.Range("B48").Value = 7688
$ws.Range("C48").Value = 403
$ws.Range("D48").Value = 1722
$ws.Range("E48").Value = 5626
$ws.Range("G48").Value = 16
$ws.Range("H48").Value = 340
